$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-25: fill in / clear individual missing-data cells ---
$ws.Range("C2").Value = 14.9
$ws.Range("F2").Value = 18.03

$ws.Range("D3").Value = $null
$ws.Range("F3").Value = $null

$ws.Range("D4").Value = -15.4

$ws.Range("D5").Value = $null

$ws.Range("C6").Value = $null

$ws.Range("D8").Value = $null
$ws.Range("F8").Value = $null

$ws.Range("C12").Value = 12.5

$ws.Range("F13").Value = 17.1

$ws.Range("C14").Value = $null

$ws.Range("D15").Value = -15.2

$ws.Range("D18").Value = -15.2

$ws.Range("D19").Value = $null
$ws.Range("F19").Value = $null

$ws.Range("C20").Value = 12.5

$ws.Range("C21").Value = 12.7

$ws.Range("D22").Value = $null

$ws.Range("C23").Value = $null
$ws.Range("D23").Value = -13.9

$ws.Range("C24").Value = $null

$ws.Range("D25").Value = -15.5
$ws.Range("F25").Value = 16.6

# --- Rows 26-35: remove the last two rows and rewrite rows 26-33 ---
$ws.Rows("34:35").Delete()

$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").Value = -20.2
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38

$ws.Range("A27").Value = "SC 101"
$ws.Range("B27").Value = $null
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = $null
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = 17

$ws.Range("A28").Value = "SC 105"
$ws.Range("B28").Value = -19.6
$ws.Range("C28").Value = 11.1
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = 17.44

$ws.Range("A29").Value = "SC 119"
$ws.Range("B29").Value = -19.5
$ws.Range("C29").Value = 11.2
$ws.Range("D29").Value = $null
$ws.Range("E29").Value = -6.8
$ws.Range("F29").Value = 18.06

$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89

$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = $null

$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").Value = $null
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53
